# Updated cryptos list on Sun Apr 30 13:42:52 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to store a literal text value (avoid Excel
# auto-converting numeric-looking strings like "326.51" or "1.010"
# into actual numbers), without leaving a stray number-format style
# behind on the cell.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "29.618.03"
Set-TextValue $ws.Range("E2") "  +0.48%  "

Set-TextValue $ws.Range("D3") "1.927.49"
Set-TextValue $ws.Range("E3") "  +0.72%  "

Set-TextValue $ws.Range("E4") "  +0.77%  "

Set-TextValue $ws.Range("D5") "326.51"
Set-TextValue $ws.Range("E5") "  +0.39%  "

Set-TextValue $ws.Range("E6") "  +0.71%  "

Set-TextValue $ws.Range("D7") "0.4825"
Set-TextValue $ws.Range("E7") "  -0.20%  "

Set-TextValue $ws.Range("D8") "0.4058"
Set-TextValue $ws.Range("E8") "  -0.33%  "

Set-TextValue $ws.Range("D9") "0.08198"
Set-TextValue $ws.Range("E9") "  +0.37%  "

Set-TextValue $ws.Range("D10") "1.008"
Set-TextValue $ws.Range("E10") "  -0.65%  "

Set-TextValue $ws.Range("D11") "23.73"
Set-TextValue $ws.Range("E11") "  +0.16%  "

Set-TextValue $ws.Range("D12") "1.925.95"
Set-TextValue $ws.Range("E12") "  +0.47%  "

Set-TextValue $ws.Range("D13") "6.068"
Set-TextValue $ws.Range("E13") "  +0.39%  "

Set-TextValue $ws.Range("D14") "7.286"
Set-TextValue $ws.Range("E14") "  +1.06%  "

Set-TextValue $ws.Range("E15") "  +0.44%  "

Set-TextValue $ws.Range("D16") "0.06860"
Set-TextValue $ws.Range("E16") "  +1.46%  "

Set-TextValue $ws.Range("E17") "  +0.64%  "

Set-TextValue $ws.Range("E18") "  -0.05%  "

Set-TextValue $ws.Range("D19") "17.59"
Set-TextValue $ws.Range("E19") "  -0.90%  "

Set-TextValue $ws.Range("D20") "1.011"
Set-TextValue $ws.Range("E20") "  +0.63%  "

Set-TextValue $ws.Range("D21") "29.606.64"
Set-TextValue $ws.Range("E21") "  +0.33%  "

Set-TextValue $ws.Range("D22") "5.654"
Set-TextValue $ws.Range("E22") "  +0.41%  "

Set-TextValue $ws.Range("D23") "11.93"
Set-TextValue $ws.Range("E23") "  +1.63%  "

Set-TextValue $ws.Range("D24") "2.198"
Set-TextValue $ws.Range("E24") "  +0.95%  "

Set-TextValue $ws.Range("D25") "2.140.65"
Set-TextValue $ws.Range("E25") "  -0.44%  "

Set-TextValue $ws.Range("D26") "156.41"
Set-TextValue $ws.Range("E26") "  -0.04%  "

Set-TextValue $ws.Range("D27") "6.360"
Set-TextValue $ws.Range("E27") "  -3.14%  "

Set-TextValue $ws.Range("D28") "19.96"
Set-TextValue $ws.Range("E28") "  -0.71%  "

Set-TextValue $ws.Range("D29") "2.086"
Set-TextValue $ws.Range("E29") "  -1.80%  "

Set-TextValue $ws.Range("D30") "120.74"
Set-TextValue $ws.Range("E30") "  +0.10%  "

Set-TextValue $ws.Range("D31") "1.002"
Set-TextValue $ws.Range("E31") "  -2.52%  "

Set-TextValue $ws.Range("D32") "0.09584"
Set-TextValue $ws.Range("E32") "  +0.41%  "

Set-TextValue $ws.Range("D33") "5.630"
Set-TextValue $ws.Range("E33") "  +1.93%  "

Set-TextValue $ws.Range("D34") "3.558"
Set-TextValue $ws.Range("E34") "  -0.08%  "

Set-TextValue $ws.Range("E35") "  -0.21%  "

Set-TextValue $ws.Range("D36") "0.06537"
Set-TextValue $ws.Range("E36") "  +6.70%  "

Set-TextValue $ws.Range("E37") "  +0.05%  "

Set-TextValue $ws.Range("D38") "1.215"
Set-TextValue $ws.Range("E38") "  +2.37%  "

Set-TextValue $ws.Range("D39") "0.5931"
Set-TextValue $ws.Range("E39") "  -0.78%  "

Set-TextValue $ws.Range("D40") "10.73"
Set-TextValue $ws.Range("E40") "  -1.73%  "

Set-TextValue $ws.Range("D41") "7.856"
Set-TextValue $ws.Range("E41") "  -2.40%  "

Set-TextValue $ws.Range("D42") "0.1844"
Set-TextValue $ws.Range("E42") "  -0.66%  "

Set-TextValue $ws.Range("D43") "2.483"
Set-TextValue $ws.Range("E43") "  +3.12%  "

Set-TextValue $ws.Range("D44") "1.242"
Set-TextValue $ws.Range("E44") "  -2.80%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D45") "12.41"
Set-TextValue $ws.Range("E45") "  -0.07%  "

$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D46") "0.07547"
Set-TextValue $ws.Range("E46") "  -0.93%  "

Set-TextValue $ws.Range("D47") "0.5550"
Set-TextValue $ws.Range("E47") "  -0.75%  "

Set-TextValue $ws.Range("D48") "1.960"
Set-TextValue $ws.Range("E48") "  +0.24%  "

Set-TextValue $ws.Range("D49") "118.26"
Set-TextValue $ws.Range("E49") "  +1.46%  "

Set-TextValue $ws.Range("D50") "2.433"
Set-TextValue $ws.Range("E50") "  +0.51%  "

Set-TextValue $ws.Range("D51") "71.99"
Set-TextValue $ws.Range("E51") "  -1.07%  "
